$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "dSF" column (column F) values for the affected rows.
# This reflects a repull/push of data with recalculated mean values.
$ws.Range("F4").Value = 7
$ws.Range("F6").Value = 10
$ws.Range("F9").Value = 7
$ws.Range("F13").Value = -4
$ws.Range("F15").Value = -3
$ws.Range("F17").Value = 7
$ws.Range("F26").Value = -3
$ws.Range("F27").Value = -6
$ws.Range("F29").Value = 3
$ws.Range("F30").Value = -5
$ws.Range("F34").Value = -4
$ws.Range("F36").Value = -6
$ws.Range("F38").Value = 4
$ws.Range("F45").Value = -3
$ws.Range("F47").Value = -1
$ws.Range("F48").Value = 6
$ws.Range("F49").Value = -1
$ws.Range("F52").Value = -2
$ws.Range("F53").Value = -7
